# Update the existing "Cards" sheet: add Art and CardType columns
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Cards")

$ws1.Range("D1").Value = "Art"
$ws1.Range("D2").Value = "half-heart"
$ws1.Range("D3").Value = "half-heart"
$ws1.Range("D4").Value = "half-heart"

$ws1.Range("E1").Value = "CardType"
$ws1.Range("E2").Value = "Skill"
$ws1.Range("E3").Value = "Skill"
$ws1.Range("E4").Value = "Skill"

# Add a new worksheet "Sheet1" right after "Cards"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

$ws2.Range("A2").Value = "Card1"
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 2
$ws2.Range("D2").Value = 4
$ws2.Range("E2").Value = 8

$ws2.Range("A3").Value = "Card2"
$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 4
$ws2.Range("D3").Value = 8
$ws2.Range("E3").Value = 16

$ws2.Range("A4").Value = "Card3"
$ws2.Range("B4").Value = 3
$ws2.Range("C4").Value = 6
$ws2.Range("D4").Value = 12
$ws2.Range("E4").Value = 24

# Set selection on the new sheet first ...
$ws2.Activate()
$ws2.Range("E16:E17").Select()

# ... then re-activate "Cards" so it remains the tab that is shown/selected
$ws1.Activate()
$ws1.Range("H7").Select()
